$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated typology mapping text (50:50 split between CDL and CDN) ---

$officesText = @'
1.8% CR/LFM+CDN/H:1/Offices
4.55% CR/LFM+CDN/H:2/Offices
10.9% CR/LFM+CDN/HBET:3-5/Offices
0.9% CR/LFM+CDN/HBET:6-/Offices
1.8% CR/LFM+CDL/H:1/Offices
4.55% CR/LFM+CDL/H:2/Offices
10.9% CR/LFM+CDL/HBET:3-5/Offices
0.9% CR/LFM+CDL/HBET:6-/Offices
3.6% MR/LWAL+CDL/H:1/Offices
7.3% MR/LWAL+CDL/H:2/Offices
25.5% MR/LWAL+CDL/HBET:3-5/Offices
17.3% S/LFM+CDL/HBET:3-5/Offices
0.9% S/LFM+CDL/HBET:6-/Offices
8.6% S/LWAL+CDL/HBET:3-5/Offices
0.5% S/LWAL+CDL/HBET:6-/Offices
0.0% S/LFM+CDL/H:1/Offices
0.0% S/LFM+CDL/H:2/Offices
0.0% S/LWAL+CDL/H:1/Offices
0.0% S/LWAL+CDL/H:2/Offices
0.0% W/LFM+CDL/H:1/Offices
0.0% W/LFM+CDL/H:2/Offices
'@

$tradeText = @'
8.2% CR/LFM+CDL/H:1/Trade
0.9% CR/LFM+CDL/H:2/Trade
8.2% CR/LFM+CDL/H:1/Trade
0.9% CR/LFM+CDL/H:2/Trade
0.0% CR/LFM+CDL/HBET:3-5/Trade
0.0% CR/LFM+CDL/HBET:6-/Trade
32.7% MR/LWAL+CDL/H:1/Trade
3.6% MR/LWAL+CDL/H:2/Trade
0.0% MR/LWAL+CDL/HBET:3-5/Trade
0.0% S/LFM+CDL/HBET:3-5/Trade
0.0% S/LFM+CDL/HBET:6-/Trade
0.0% S/LWAL+CDL/HBET:3-5/Trade
0.0% S/LWAL+CDL/HBET:6-/Trade
17.8% S/LFM+CDL/H:1/Trade
0.4% S/LFM+CDL/H:2/Trade
8.9% S/LWAL+CDL/H:1/Trade
0.2% S/LWAL+CDL/H:2/Trade
17.8% W/LFM+CDL/H:1/Trade
0.4% W/LFM+CDL/H:2/Trade
'@

$hotelsText = @'
2.2% CR/LFM+CDN/H:1/Hotels
0.9% CR/LFM+CDN/H:2/Hotels
5.3% CR/LFM+CDN/HBET:3-5/Hotels
0.45% CR/LFM+CDN/HBET:6-/Hotels
2.2% CR/LFM+CDL/H:1/Hotels
0.9% CR/LFM+CDL/H:2/Hotels
5.3% CR/LFM+CDL/HBET:3-5/Hotels
0.45% CR/LFM+CDL/HBET:6-/Hotels
 3.5% MR/LWAL+CDL/H:1/Hotels
 7.1% MR/LWAL+CDL/H:2/Hotels
 25.1% MR/LWAL+CDL/HBET:3-5/Hotels
 17.2% S/LFM+CDL/HBET:3-5/Hotels
 2.9% S/LFM+CDL/HBET:6-/Hotels
 8.4% S/LWAL+CDL/HBET:3-5/Hotels
 0.4% S/LWAL+CDL/HBET:6-/Hotels
 0.0% S/LFM+CDL/H:1/Hotels
 0.0% S/LFM+CDL/H:2/Hotels
 0.0% S/LWAL+CDL/H:1/Hotels
 0.0% S/LWAL+CDL/H:2/Hotels
 14.2% W/LFM+CDL/H:1/Hotels
 3.5% W/LFM+CDL/H:2/Hotels
'@

$ws.Range("B2").Value2 = $officesText
$ws.Range("C2").Value2 = $tradeText
$ws.Range("D2").Value2 = $hotelsText

# --- Wrap text for the long mapping cells (new style entry) ---
$ws.Range("B2:D2").WrapText = $true

# --- Column widths (target ~32.5 / 32.33 / 32.5 "characters") ---
$ws.Columns.Item(2).ColumnWidth = 190/6
$ws.Columns.Item(3).ColumnWidth = 189/6
$ws.Columns.Item(4).ColumnWidth = 190/6

# --- Row height for the wrapped row ---
$ws.Rows.Item(2).RowHeight = 335

# --- Update selection to match authored state ---
[void]$ws.Range("F2").Select()
